$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns remain stored as text,
# matching the source data (values like "23.257.42" or "0.9989" would
# otherwise be auto-coerced into numbers by Excel).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "23.257.42"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "1.603.47"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "0.9992"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "303.23"
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("D7").Value = "0.3777"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "51.72"
$ws.Range("E8").Value = "  +3.46%  "
$ws.Range("D9").Value = "0.3635"
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("D10").Value = "1.272"
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("D11").Value = "0.08136"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").Value = "0.9992"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").Value = "6.604"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "7.421"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").Value = "0.00001248"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "1.604.57"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "93.90"
$ws.Range("E18").Value = "  +1.93%  "
$ws.Range("D19").Value = "0.06888"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "18.14"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").Value = "6.548"
$ws.Range("E21").Value = "  -0.61%  "
$ws.Range("D22").Value = "0.9992"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "12.97"
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("D24").Value = "23.262.71"
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "2.392"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "3.010"
$ws.Range("E26").Value = "  +7.51%  "
$ws.Range("D27").Value = "21.26"
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("D28").Value = "149.98"
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("D29").Value = "5.255"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").Value = "134.15"
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("D31").Value = "2.375"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("D32").Value = "6.779"
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("D33").Value = "1.780.72"
$ws.Range("E33").Value = "  +0.35%  "
$ws.Range("D34").Value = "0.9678"
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("D35").Value = "0.07528"
$ws.Range("E35").Value = "  -2.24%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "10.28"
$ws.Range("E36").Value = "  -1.76%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.02734"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").Value = "0.2534"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").Value = "0.08810"
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("D40").Value = "6.091"
$ws.Range("E40").Value = "  -3.20%  "
$ws.Range("D41").Value = "1.373"
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("D42").Value = "0.7123"
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("D43").Value = "12.54"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("D44").Value = "15.63"
$ws.Range("E44").Value = "  +2.12%  "
$ws.Range("D45").Value = "0.6563"
$ws.Range("E45").Value = "  -1.17%  "
$ws.Range("D46").Value = "2.320"
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("D47").Value = "4.015"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").Value = "132.37"
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").Value = "0.07965"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").Value = "1.205"
$ws.Range("E50").Value = "  -3.29%  "
$ws.Range("D51").Value = "1.207"
$ws.Range("E51").Value = "  +0.20%  "
